$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: Farhan ---
$ws.Range("B2").Value = "Farhan"
$ws.Range("C2").NumberFormat = "@"
$ws.Range("C2").Value = "+60177496576"
$ws.Range("D2").Value = "farhan257@gmail.com"
$ws.Range("G2").Value = "[{'Country': '', 'State': 'Kuala Lumpur', 'City': ''}]"
$ws.Range("H2").Value = "[{'field_of_study': 'Software Engineering', 'level': 'Bachelor', 'cgpa': '3.54', 'university': 'Universiti Teknologi Malaysia (UTM)', 'start_date': 'N/A', 'year_of_graduation': 'N/A'}, {'field_of_study': 'Software Engineering', 'level': 'Master of Software Engineering', 'cgpa': '3.98', 'university': 'Universiti Teknologi Malaysia (UTM)', 'start_date': 'N/A', 'year_of_graduation': 'N/A'}]"
$ws.Range("I2").Value = "[]"
$ws.Range("J2").Value = "['Agile Software Development', 'Requirement analysis', 'DevOps', 'Linux/Unix environment', 'CI/CD pipeline (Jenkins)', 'Containerization (Kubernetes & Docker)', 'Git (Git Bash, GitHub, GitLab)', 'Software Architecture & Design', 'Software Testing (Cucumber, JMeter)', 'Virtualization (VMWare/VirtualBox)', 'Software Documentation', 'Software Standards (ISO/IEC/IEEE)', 'SQL Server Management', 'Laravel', 'Azure', 'Programming Languages (Java, PHP, C++, SQL, JavaScript, Python)', 'Web Development', 'Microsoft Office', 'Enterprise Architect', 'Epicor', 'Visual Studio, Android Studio, Eclipse, Spring Tool Suite, Ionic', 'Wireshark/Tshark']"
$ws.Range("K2").Value = "[]"
$ws.Range("L2").Value = "['English', 'Bahasa Malaysia']"
$ws.Range("M2").Value = "[{'job_title': 'Consultant', 'job_company': 'Finsoft Consulting Sdn Bhd', 'Industries': 'N/A', 'start_date': '2022-04-01', 'end_date': '2022-07-01', 'job_location': 'N/A'}, {'job_title': 'Software Engineer', 'job_company': 'Axacute', 'Industries': 'N/A', 'start_date': '2020-07-01', 'end_date': '2021-10-01', 'job_location': 'N/A'}, {'job_title': 'Intern', 'job_company': 'Openet', 'Industries': 'N/A', 'start_date': '2019-07-01', 'end_date': '2019-12-01', 'job_location': 'N/A'}]"

# --- Row 3: MOHAMAD AMIR AFIFIE ---
$ws.Range("B3").Value = "MOHAMAD AMIR AFIFIE"
$ws.Range("C3").Value = "0111 - 488 3732"
$ws.Range("D3").Value = "amirafifie@gmail.com"
$ws.Range("E3").Value = "N/A"
$ws.Range("G3").Value = "['Malaysia', 'Cyberjaya', '']"
$ws.Range("H3").Value = "[{'field_of_study': 'Computer Science', 'level': ""Bachelor's Degree"", 'cgpa': 'N/A', 'university': 'NATIONAL UNIVERSITY OF MALAYSIA', 'start_date': '2017', 'year_of_graduation': '2021'}, {'field_of_study': 'Science', 'level': 'Foundation', 'cgpa': 'N/A', 'university': 'UNIVERSITI TEKNOLOGI MARA', 'start_date': '2016', 'year_of_graduation': '2017'}]"
$ws.Range("I3").Value = "[]"
$ws.Range("J3").Value = "['Data Visualisation', 'Machine Learning', 'Data Analysis', 'Python', 'Java', 'SQL', 'C', 'VBA', 'Power BI', 'Excel', 'Scikit-learn', 'NLTK', 'Vader', 'Textblob', 'Pandas', 'Matplotlib', 'Numpy', 'Problem-Solving', 'Team Player', 'Communication']"
$ws.Range("K3").Value = "[]"
$ws.Range("L3").Value = "[]"
$ws.Range("M3").Value = "[{'job_title': 'HR DATA ANALYST INTERN', 'job_company': 'Safran Landing System', 'Industries': 'N/A', 'start_date': '2022-08-01', 'end_date': '2023-01-31', 'job_location': 'N/A', 'job_duration': 0.42}]"

Write-Output "Updated rows 2 and 3 successfully."
